$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '87.319.74'
$ws.Range('E2').Value = '  +3.05%  '

$ws.Range('D3').Value = '3.242.08'
$ws.Range('E3').Value = '  -1.91%  '

$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5').Value = "'210.50"
$ws.Range('E5').Value = '  -4.22%  '

$ws.Range('D6').Value = "'624.04"
$ws.Range('E6').Value = '  -1.98%  '

$ws.Range('D7').Value = "'0.374"
$ws.Range('E7').Value = '  +16.09%  '

$ws.Range('D8').Value = "'0.683"
$ws.Range('E8').Value = '  +15.20%  '

$ws.Range('D9').Value = "'1.00"
$ws.Range('E9').Value = '  +0.10%  '

$ws.Range('D10').Value = '3.241.33'
$ws.Range('E10').Value = '  -2.66%  '

$ws.Range('D11').Value = "'0.572"
$ws.Range('E11').Value = '  -4.45%  '

$ws.Range('D12').Value = "'0.184"
$ws.Range('E12').Value = '  +10.90%  '

$ws.Range('D13').Value = "'0.0000258"
$ws.Range('E13').Value = '  -6.97%  '

$ws.Range('D14').Value = '3.850.90'
$ws.Range('E14').Value = '  -1.66%  '

$ws.Range('D15').Value = "'33.78"
$ws.Range('E15').Value = '  -1.04%  '

$ws.Range('D16').Value = "'5.31"
$ws.Range('E16').Value = '  -2.44%  '

$ws.Range('D17').Value = '87.271.35'
$ws.Range('E17').Value = '  +3.27%  '

$ws.Range('D18').Value = '3.256.18'
$ws.Range('E18').Value = '  -1.41%  '

$ws.Range('D19').Value = "'3.11"
$ws.Range('E19').Value = '  -3.07%  '

$ws.Range('D20').Value = "'13.94"
$ws.Range('E20').Value = '  -4.89%  '

$ws.Range('D21').Value = "'431.12"
$ws.Range('E21').Value = '  -2.15%  '

$ws.Range('D22').Value = "'8.82"
$ws.Range('E22').Value = '  -4.54%  '

$ws.Range('E23').Value = '  +1.76%  '

$ws.Range('D24').Value = "'7.33"
$ws.Range('E24').Value = '  -0.53%  '

$ws.Range('D25').Value = "'12.35"
$ws.Range('E25').Value = '  +1.32%  '

$ws.Range('D26').Value = "'5.09"
$ws.Range('E26').Value = '  -7.38%  '

$ws.Range('D27').Value = '3.373.84'
$ws.Range('E27').Value = '  -2.51%  '

$ws.Range('D28').Value = "'76.08"
$ws.Range('E28').Value = '  -2.60%  '

$ws.Range('D29').Value = "'0.0000129"
$ws.Range('E29').Value = '  -2.02%  '

$ws.Range('D30').Value = "'1.00"
$ws.Range('E30').Value = '  +0.01%  '

$ws.Range('D31').Value = "'0.183"
$ws.Range('E31').Value = '  +9.96%  '

$ws.Range('D32').Value = "'0.998"
$ws.Range('E32').Value = '  -0.30%  '

$ws.Range('D33').Value = "'8.72"
$ws.Range('E33').Value = '  -6.00%  '

$ws.Range('D34').Value = "'552.64"
$ws.Range('E34').Value = '  -7.62%  '

$ws.Range('E35').Value = '  -12.07%  '

$ws.Range('E36').Value = '  -4.76%  '

$ws.Range('D37').Value = "'6.93"

$ws.Range('E38').Value = '  -10.84%  '

$ws.Range('D39').Value = "'22.36"
$ws.Range('E39').Value = '  -3.98%  '

$ws.Range('D40').Value = "'1.00"
$ws.Range('E40').Value = '  +0.14%  '

$ws.Range('D41').Value = "'21.73"
$ws.Range('E41').Value = '  +3.74%  '

$ws.Range('D42').Value = "'0.391"
$ws.Range('E42').Value = '  -6.30%  '

$ws.Range('D43').Value = "'1.99"
$ws.Range('E43').Value = '  -2.72%  '

$ws.Range('D44').Value = "'2.92"
$ws.Range('E44').Value = '  -6.30%  '

$ws.Range('E45').Value = '  -0.01%  '

$ws.Range('D46').Value = "'151.55"
$ws.Range('E46').Value = '  -4.85%  '

$ws.Range('D47').Value = "'178.11"
$ws.Range('E47').Value = '  -6.57%  '

$ws.Range('D48').Value = "'44.50"
$ws.Range('E48').Value = '  -1.42%  '

$ws.Range('D49').Value = "'1.28"
$ws.Range('E49').Value = '  -5.23%  '

$ws.Range('D50').Value = "'4.20"
$ws.Range('E50').Value = '  -1.13%  '

$ws.Range('B51').Value = 'Hedera'
$ws.Range('C51').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D51').Value = "'0.0662"
$ws.Range('E51').Value = '  +18.38%  '
